# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the first data row on both the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: column E = Correspond Handoff Datetime, column H = Correspond Handback DateTime
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-21 11:01:06"
$wsZh.Range("H2").Value = "2016-03-21 11:01:46"

# de-de sheet: column E = Correspond Handoff Datetime, column H = Correspond Handback DateTime
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-21 11:01:16"
$wsDe.Range("H2").Value = "2016-03-21 11:01:52"
